$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Helper: assign a text value to a cell while preserving its original "General"
# (no explicit style) formatting. Some replacement values look like plain numbers
# (e.g. "252.65"); without forcing a temporary Text format Excel would silently
# convert them to numeric values, which would not match the source data (plain
# text cells). We flip the format to Text only long enough to assign the value,
# then restore the cell style so no stray formatting is left behind.
function Set-TextValue($cell, $value) {
    $rng = $ws.Range($cell)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

$ws.Range('D2').Value = '36.542.69'
$ws.Range('E2').Value = '  -0.36%  '
$ws.Range('D3').Value = '2.098.95'
$ws.Range('E3').Value = '  +9.56%  '
$ws.Range('E4').Value = '  +0.02%  '
Set-TextValue 'D5' '252.65'
$ws.Range('E5').Value = '  +0.81%  '
Set-TextValue 'D6' '0.657'
$ws.Range('E6').Value = '  -6.36%  '
$ws.Range('E7').Value = '  +0.01%  '
Set-TextValue 'D8' '47.36'
$ws.Range('E8').Value = '  +5.46%  '
Set-TextValue 'D9' '60.22'
$ws.Range('E9').Value = '  +3.20%  '
Set-TextValue 'D10' '0.377'
$ws.Range('E10').Value = '  +1.46%  '
Set-TextValue 'D11' '0.0743'
$ws.Range('E11').Value = '  -2.69%  '
$ws.Range('E12').Value = '  -0.16%  '
Set-TextValue 'D13' '14.66'
$ws.Range('E13').Value = '  +0.77%  '
$ws.Range('D14').Value = '2.404.24'
$ws.Range('E14').Value = '  +9.48%  '
Set-TextValue 'D15' '0.825'
$ws.Range('E15').Value = '  +0.96%  '
$ws.Range('D16').Value = '2.102.01'
$ws.Range('E16').Value = '  +9.58%  '
Set-TextValue 'D17' '5.08'
$ws.Range('E17').Value = '  -1.08%  '
$ws.Range('D18').Value = '36.501.28'
$ws.Range('E18').Value = '  -0.45%  '
Set-TextValue 'D19' '72.70'
$ws.Range('E19').Value = '  -2.62%  '
$ws.Range('D20').Value = '0.0₃0828'
$ws.Range('E20').Value = '  -3.97%  '
Set-TextValue 'D21' '13.15'
$ws.Range('E21').Value = '  -1.77%  '
Set-TextValue 'D22' '239.64'
$ws.Range('E22').Value = '  -4.33%  '
$ws.Range('E23').Value = '  -1.17%  '
$ws.Range('E24').Value = '  +0.04%  '
Set-TextValue 'D25' '2.46'
$ws.Range('E25').Value = '  -6.74%  '
Set-TextValue 'D26' '170.76'
$ws.Range('E26').Value = '  +1.20%  '
Set-TextValue 'D27' '21.39'
$ws.Range('E27').Value = '  +14.19%  '
Set-TextValue 'D28' '9.12'
$ws.Range('E28').Value = '  +3.76%  '
Set-TextValue 'D29' '1.99'
$ws.Range('E29').Value = '  -9.94%  '
Set-TextValue 'D30' '28.80'
$ws.Range('E30').Value = '  +62.13%  '
$ws.Range('E31').Value = '  -5.29%  '
Set-TextValue 'D32' '4.45'
$ws.Range('E32').Value = '  -2.83%  '
Set-TextValue 'D33' '0.0613'
$ws.Range('E33').Value = '  -1.06%  '
Set-TextValue 'D34' '0.0924'
$ws.Range('E34').Value = '  +2.72%  '
Set-TextValue 'D35' '0.989'
$ws.Range('E35').Value = '  +12.35%  '
Set-TextValue 'D36' '2.41'
$ws.Range('E36').Value = '  +19.28%  '
$ws.Range('E37').Value = '  +0.05%  '
$ws.Range('E38').Value = '  -2.83%  '
Set-TextValue 'D39' '4.08'
$ws.Range('E39').Value = '  -5.95%  '
$ws.Range('E40').Value = '  -11.35%  '
$ws.Range('E41').Value = '  +5.63%  '
$ws.Range('E42').Value = '  -2.16%  '
Set-TextValue 'D43' '97.58'
$ws.Range('E43').Value = '  -8.03%  '
Set-TextValue 'D44' '2.75'
$ws.Range('E44').Value = '  -5.98%  '
Set-TextValue 'D45' '15.93'
$ws.Range('E45').Value = '  -9.08%  '
$ws.Range('D46').Value = '1.324.78'
$ws.Range('E46').Value = '  -1.53%  '
Set-TextValue 'D47' '0.0842'
$ws.Range('E47').Value = '  +3.58%  '
Set-TextValue 'D48' '6.99'
$ws.Range('E48').Value = '  +8.63%  '
$ws.Range('B49').Value = 'MXToken'
$ws.Range('C49').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
Set-TextValue 'D49' '2.85'
$ws.Range('E49').Value = '  +1.53%  '
$ws.Range('B50').Value = 'RocketPoolETH'
$ws.Range('C50').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D50').Value = '2.294.34'
$ws.Range('E50').Value = '  +9.63%  '
$ws.Range('E51').Value = '  -6.29%  '
